$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.553.78"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").Value = "2.566.41"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.75"
$ws.Range("E5").Value = "  -0.84%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.54"
$ws.Range("E6").Value = "  -3.25%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("E9").Value = "  -2.72%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.55"
$ws.Range("E10").Value = "  -0.78%  "

$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.351"
$ws.Range("E12").Value = "  -1.65%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "26.82"
$ws.Range("E13").Value = "  -4.20%  "

$ws.Range("D14").Value = "3.026.53"
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("D15").Value = "62.546.44"
$ws.Range("E15").Value = "  -0.83%  "

$ws.Range("E16").Value = "  -2.56%  "

$ws.Range("D17").Value = "2.569.80"
$ws.Range("E17").Value = "  +0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.14"
$ws.Range("E18").Value = "  -2.73%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.99"
$ws.Range("E19").Value = "  -1.08%  "

$ws.Range("E20").Value = "  -1.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.65"
$ws.Range("E21").Value = "  -2.98%  "

$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.15"
$ws.Range("E23").Value = "  +1.67%  "

$ws.Range("E24").Value = "  -3.73%  "

$ws.Range("E25").Value = "  -4.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.49"
$ws.Range("E26").Value = "  +0.27%  "

$ws.Range("E27").Value = "  -0.22%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.94"
$ws.Range("E28").Value = "  -3.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.19"
$ws.Range("E29").Value = "  -4.06%  "

$ws.Range("E30").Value = "  -2.17%  "

$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "457.34"
$ws.Range("E31").Value = "  +4.44%  "

$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").Value = "0.0₃0799"
$ws.Range("E32").Value = "  -3.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "176.57"
$ws.Range("E33").Value = "  -0.58%  "

$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.395"
$ws.Range("E36").Value = "  -3.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.83"
$ws.Range("E37").Value = "  -2.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.45"
$ws.Range("E38").Value = "  -2.21%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("E40").Value = "  -4.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "158.10"
$ws.Range("E41").Value = "  +3.82%  "

$ws.Range("E42").Value = "  -3.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.628"
$ws.Range("E43").Value = "  +3.09%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "20.82"
$ws.Range("E44").Value = "  -2.55%  "

$ws.Range("E45").Value = "  -4.38%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0959"
$ws.Range("E46").Value = "  -1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0233"
$ws.Range("E47").Value = "  -3.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.96"
$ws.Range("E48").Value = "  -2.75%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "11.40"
$ws.Range("E49").Value = "  +0.17%  "

$ws.Range("E50").Value = "  -4.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.956"
$ws.Range("E51").Value = "  +3.25%  "
